# Regenerate the Sheet1 report for "config 03": the underlying Python run
# failed with "can't convert type 'dict' to numerator/denominator" for the
# per-method error-rate metrics, and the before/after-fix index stats
# collapsed to zero. The after_fix_* columns (J:L) are no longer produced,
# so drop them and shift nothing else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused after_fix_mean / after_fix_variance / after_fix_std
# columns (J:L) entirely, shrinking the used range from A1:L17 to A1:I17.
$ws.Range("J:L").EntireColumn.Delete()

# --- Header row -----------------------------------------------------------
$ws.Range("C1").Value = "Erro"
$ws.Range("D1").Value = "before_fix_mean"
$ws.Range("E1").Value = "before_fix_variance"
$ws.Range("F1").Value = "before_fix_std"
$ws.Range("G1").Value = "after_fix_mean"
$ws.Range("H1").Value = "after_fix_variance"
$ws.Range("I1").Value = "after_fix_std"

$errMsg = "can't convert type 'dict' to numerator/denominator"

# --- Per-dataset metric rows (error rate computation now fails) -----------
$metricRows = 2, 3, 4, 5, 10, 11, 12, 13
foreach ($r in $metricRows) {
    $ws.Range("C${r}").Value = $errMsg
    $ws.Range("D${r}:I${r}").ClearContents()
}

# --- "before_fix" index-stat rows: mean/variance/std now 0, and the ------
# stats have moved from columns G:I to D:F in the new layout.
$beforeFixRows = 6, 7, 14, 15
foreach ($r in $beforeFixRows) {
    $ws.Range("D${r}").Value = 0
    $ws.Range("E${r}").Value = 0
    $ws.Range("F${r}").Value = 0
    $ws.Range("G${r}:I${r}").ClearContents()
}

# --- "after_fix" index-stat rows: mean/variance/std now 0, staying in ----
# columns G:I (the old J:K:L data, now gone).
$afterFixRows = 8, 9, 16, 17
foreach ($r in $afterFixRows) {
    $ws.Range("G${r}").Value = 0
    $ws.Range("H${r}").Value = 0
    $ws.Range("I${r}").Value = 0
}
